{"js": "// Replace the three-digit-by-one-digit multiplication answers in the\n// practice table with the new values from the commit.\nconst replacements = [\n  [\"958\u00d76=5748\", \"414\u00d74=1656\"],\n  [\"975\u00d72=1950\", \"575\u00d79=5175\"],\n  [\"830\u00d72=1660\", \"620\u00d72=1240\"],\n  [\"976\u00d79=8784\", \"882\u00d79=7938\"],\n  [\"162\u00d77=1134\", \"543\u00d74=2172\"],\n  [\"781\u00d78=6248\", \"454\u00d72=908\"],\n  [\"697\u00d77=4879\", \"584\u00d72=1168\"],\n  [\"881\u00d78=7048\", \"170\u00d77=1190\"],\n  [\"847\u00d72=1694\", \"641\u00d75=3205\"],\n  [\"677\u00d73=2031\", \"952\u00d73=2856\"],\n  [\"864\u00d74=3456\", \"112\u00d78=896\"],\n  [\"490\u00d74=1960\", \"803\u00d75=4015\"],\n  [\"988\u00d79=8892\", \"689\u00d72=1378\"],\n  [\"934\u00d74=3736\", \"411\u00d75=2055\"],\n  [\"996\u00d78=7968\", \"724\u00d74=2896\"],\n  [\"361\u00d72=722\", \"815\u00d77=5705\"],\n  [\"770\u00d76=4620\", \"311\u00d78=2488\"],\n  [\"326\u00d75=1630\", \"405\u00d79=3645\"],\n  [\"808\u00d74=3232\", \"535\u00d78=4280\"],\n  [\"111\u00d77=777\", \"287\u00d73=861\"],\n  [\"348\u00d75=1740\", \"360\u00d77=2520\"],\n  [\"879\u00d72=1758\", \"132\u00d72=264\"],\n  [\"880\u00d78=7040\", \"852\u00d76=5112\"],\n  [\"548\u00d76=3288\", \"284\u00d75=1420\"],\n  [\"676\u00d77=4732\", \"899\u00d79=8091\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication answers in the\n# practice table with the new values from the commit.\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n    @(\"958\u00d76=5748\", \"414\u00d74=1656\"),\n    @(\"975\u00d72=1950\", \"575\u00d79=5175\"),\n    @(\"830\u00d72=1660\", \"620\u00d72=1240\"),\n    @(\"976\u00d79=8784\", \"882\u00d79=7938\"),\n    @(\"162\u00d77=1134\", \"543\u00d74=2172\"),\n    @(\"781\u00d78=6248\", \"454\u00d72=908\"),\n    @(\"697\u00d77=4879\", \"584\u00d72=1168\"),\n    @(\"881\u00d78=7048\", \"170\u00d77=1190\"),\n    @(\"847\u00d72=1694\", \"641\u00d75=3205\"),\n    @(\"677\u00d73=2031\", \"952\u00d73=2856\"),\n    @(\"864\u00d74=3456\", \"112\u00d78=896\"),\n    @(\"490\u00d74=1960\", \"803\u00d75=4015\"),\n    @(\"988\u00d79=8892\", \"689\u00d72=1378\"),\n    @(\"934\u00d74=3736\", \"411\u00d75=2055\"),\n    @(\"996\u00d78=7968\", \"724\u00d74=2896\"),\n    @(\"361\u00d72=722\", \"815\u00d77=5705\"),\n    @(\"770\u00d76=4620\", \"311\u00d78=2488\"),\n    @(\"326\u00d75=1630\", \"405\u00d79=3645\"),\n    @(\"808\u00d74=3232\", \"535\u00d78=4280\"),\n    @(\"111\u00d77=777\", \"287\u00d73=861\"),\n    @(\"348\u00d75=1740\", \"360\u00d77=2520\"),\n    @(\"879\u00d72=1758\", \"132\u00d72=264\"),\n    @(\"880\u00d78=7040\", \"852\u00d76=5112\"),\n    @(\"548\u00d76=3288\", \"284\u00d75=1420\"),\n    @(\"676\u00d77=4732\", \"899\u00d79=8091\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $true, $newText, $wdReplaceAll)\n}\n"}
